# Applies the cell value updates to the Balmung_Profits workbook sheets
# as described by the source commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1145.7084
$ws.Range("I41").Value = 672.4545000000001
$ws.Range("J41").Value = 1546.1538
$ws.Range("K41").Value = 672.4545000000001
$ws.Range("L41").Value = 1546.1538
$ws.Range("M41").Value = -232.4545000000001
$ws.Range("N41").Value = -2426.1538
$ws.Range("H55").Value = 306.57144
$ws.Range("I55").Value = 306.57144
$ws.Range("K55").Value = 306.57144
$ws.Range("M55").Value = -92.57144
$ws.Range("H58").Value = 505.8
$ws.Range("J58").Value = 498.5
$ws.Range("L58").Value = 1495.5
$ws.Range("N58").Value = -1795.5
$ws.Range("H98").Value = 4565.1724
$ws.Range("I98").Value = 3732
$ws.Range("K98").Value = 3732
$ws.Range("M98").Value = -2234
$ws.Range("H122").Value = 4565.1724
$ws.Range("I122").Value = 3732
$ws.Range("K122").Value = 11196
$ws.Range("M122").Value = -8746
$ws.Range("H132").Value = 1760
$ws.Range("I132").Value = 1636.15
$ws.Range("K132").Value = 4908.450000000001
$ws.Range("M132").Value = -2378.450000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5002783
$ws.Range("I39").Value = 5002783
$ws.Range("K39").Value = 5002783
$ws.Range("M39").Value = -5002263
$ws.Range("H74").Value = 931881.1
$ws.Range("I74").Value = 1654.3334
$ws.Range("J74").Value = 1490017.2
$ws.Range("K74").Value = 1654.3334
$ws.Range("L74").Value = 1490017.2
$ws.Range("M74").Value = -780.3334
$ws.Range("N74").Value = -1491765.2
$ws.Range("H77").Value = 931881.1
$ws.Range("I77").Value = 1654.3334
$ws.Range("J77").Value = 1490017.2
$ws.Range("K77").Value = 8271.666999999999
$ws.Range("L77").Value = 7450086
$ws.Range("M77").Value = -3903.666999999999
$ws.Range("N77").Value = -7458822
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 3424.4773
$ws.Range("I132").Value = 3424.4773
$ws.Range("K132").Value = 10273.4319
$ws.Range("M132").Value = -7743.4319

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 50000
$ws.Range("J57").Value = 50000
$ws.Range("L57").Value = 50000
$ws.Range("N57").Value = -51440
$ws.Range("H58").Value = 37500
$ws.Range("J58").Value = 37500
$ws.Range("L58").Value = 37500
$ws.Range("N58").Value = -38088
$ws.Range("H81").Value = 31394.143
$ws.Range("J81").Value = 31394.143
$ws.Range("L81").Value = 31394.143
$ws.Range("N81").Value = -33516.143
$ws.Range("H84").Value = 31394.143
$ws.Range("J84").Value = 31394.143
$ws.Range("L84").Value = 94182.429
$ws.Range("N84").Value = -104790.429
$ws.Range("H105").Value = 14760.556
$ws.Range("I105").Value = 15906.571
$ws.Range("K105").Value = 15906.571
$ws.Range("M105").Value = -14159.571
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H134").Value = 20456664
$ws.Range("I134").Value = 2072.1
$ws.Range("K134").Value = 6216.299999999999
$ws.Range("M134").Value = -3681.299999999999
$ws.Range("H135").Value = 60553.668
$ws.Range("J135").Value = 60553.668
$ws.Range("L135").Value = 60553.668
$ws.Range("N135").Value = -70693.66800000001
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H140").Value = 70995.5
$ws.Range("J140").Value = 70995.5
$ws.Range("L140").Value = 70995.5
$ws.Range("N140").Value = -81355.5
$ws.Range("H141").Value = 73984
$ws.Range("J141").Value = 73984
$ws.Range("L141").Value = 73984
$ws.Range("N141").Value = -84344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3297
$ws.Range("I122").Value = 3546.7
$ws.Range("K122").Value = 10640.1
$ws.Range("M122").Value = -8190.099999999999
$ws.Range("H132").Value = 18320.28
$ws.Range("I132").Value = 21573.24
$ws.Range("J132").Value = 3534.0908
$ws.Range("K132").Value = 64719.72
$ws.Range("L132").Value = 10602.2724
$ws.Range("M132").Value = -62189.72
$ws.Range("N132").Value = -15662.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 246.41379
$ws.Range("I2").Value = 200.8
$ws.Range("J2").Value = 347.77777
$ws.Range("K2").Value = 1204.8
$ws.Range("L2").Value = 2086.66662
$ws.Range("M2").Value = -1091.8
$ws.Range("N2").Value = -2312.66662
$ws.Range("H4").Value = 1000367.7
$ws.Range("I4").Value = 1148457
$ws.Range("J4").Value = 333965.66
$ws.Range("K4").Value = 3445371
$ws.Range("L4").Value = 1001896.98
$ws.Range("M4").Value = -3445259
$ws.Range("N4").Value = -1002120.98
$ws.Range("H69").Value = 8495.049999999999
$ws.Range("J69").Value = 6406.5293
$ws.Range("L69").Value = 19219.5879
$ws.Range("N69").Value = -20841.5879
$ws.Range("H72").Value = 8495.049999999999
$ws.Range("J72").Value = 6406.5293
$ws.Range("L72").Value = 57658.7637
$ws.Range("N72").Value = -65770.76370000001
$ws.Range("H97").Value = 376.44446
$ws.Range("I97").Value = 405.375
$ws.Range("K97").Value = 1216.125
$ws.Range("M97").Value = -720.125
$ws.Range("H121").Value = 864.875
$ws.Range("I121").Value = 624.6
$ws.Range("J121").Value = 1265.3334
$ws.Range("K121").Value = 1873.8
$ws.Range("L121").Value = 3796.0002
$ws.Range("M121").Value = -563.8000000000002
$ws.Range("N121").Value = -6416.0002
$ws.Range("H131").Value = 5053232
$ws.Range("I131").Value = 9092638
$ws.Range("K131").Value = 27277914
$ws.Range("M131").Value = -27272874

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2005.8684
$ws.Range("I122").Value = 1755.4849
$ws.Range("J122").Value = 3658.4
$ws.Range("K122").Value = 5266.4547
$ws.Range("L122").Value = 10975.2
$ws.Range("M122").Value = -2816.4547
$ws.Range("N122").Value = -15875.2
$ws.Range("H132").Value = 574122.4399999999
$ws.Range("I132").Value = 5404.0356
$ws.Range("J132").Value = 1510835.1
$ws.Range("K132").Value = 16212.1068
$ws.Range("L132").Value = 4532505.300000001
$ws.Range("M132").Value = -13682.1068
$ws.Range("N132").Value = -4537565.300000001
$ws.Range("H140").Value = 106720.336
$ws.Range("J140").Value = 110096.75
$ws.Range("L140").Value = 110096.75
$ws.Range("N140").Value = -120456.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2696.8147
$ws.Range("I61").Value = 2646.182
$ws.Range("J61").Value = 2919.6
$ws.Range("K61").Value = 2646.182
$ws.Range("L61").Value = 2919.6
$ws.Range("M61").Value = -2444.182
$ws.Range("N61").Value = -3323.6
$ws.Range("H113").Value = 2696.8147
$ws.Range("I113").Value = 2646.182
$ws.Range("J113").Value = 2919.6
$ws.Range("K113").Value = 2646.182
$ws.Range("L113").Value = 2919.6
$ws.Range("M113").Value = -476.1819999999998
$ws.Range("N113").Value = -7259.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20001
$ws.Range("I2").Value = 30002
$ws.Range("K2").Value = 30002
$ws.Range("M2").Value = -29890
$ws.Range("H107").Value = 1021078.8
$ws.Range("I107").Value = 699.1667
$ws.Range("J107").Value = 7143357
$ws.Range("K107").Value = 2097.5001
$ws.Range("L107").Value = 21430071
$ws.Range("M107").Value = -177.5001000000002
$ws.Range("N107").Value = -21433911
$ws.Range("H136").Value = 42325.56
$ws.Range("I136").Value = 56470.332
$ws.Range("J136").Value = 5953.2856
$ws.Range("K136").Value = 169410.996
$ws.Range("L136").Value = 17859.8568
$ws.Range("M136").Value = -166860.996
$ws.Range("N136").Value = -22959.8568

